# Exhibit 1 trade-data refresh: extend series through April, shift "(R)" flag
# from February to April, and update the "last month published" snapshot row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels / narrative text -----------------------------------------
# Quarterly sub-header now spans Jan-May instead of Jan-Mar (one more month
# of data has been added to each year block).
foreach ($r in @(9, 24, 38)) {
    $ws.Cells.Item($r, 1).Value = "Jan. - May"
}

# February is no longer the most-recently-revised month; April is. Row 40
# drops the "(R)" suffix, row 42 (April) gains it, and row 42/43 (April/May)
# get real data instead of being blank placeholders.
$ws.Cells.Item(40, 1).Value = "February"
$ws.Cells.Item(42, 1).Value = "April (R)"

# "Last month published" note/snapshot now reflects April, not February.
$ws.Cells.Item(51, 1).Value = "April data as published last month:"

# --- Numeric data refresh -------------------------------------------------
# Row number -> new values for columns B..J (Balance/Exports/Imports x Total/Goods/Services).
$updates = @{
    8 = @(-576341, -861515, 285174, 2528367, 1652072, 876295, 3104708, 2513587, 591121);
    9 = @(-248605, -367685, 119080, 1057295, 694430, 362865, 1305900, 1062115, 243785);
    10 = @(-48818, -72422, 23604, 209087, 137716, 71371, 257905, 210138, 47768);
    11 = @(-48032, -71831, 23799, 210133, 138255, 71878, 258165, 210086, 48079);
    12 = @(-49777, -73355, 23578, 213813, 141183, 72630, 263590, 214538, 49052);
    13 = @(-50074, -73654, 23580, 210289, 137284, 73005, 260363, 210937, 49425);
    14 = @(-51904, -76423, 24520, 213973, 139993, 73980, 265877, 216416, 49461);
    15 = @(-50390, -74285, 23895, 210575, 136744, 73831, 260965, 211029, 49936);
    16 = @(-49959, -73237, 23277, 211469, 138532, 72937, 261428, 211768, 49660);
    17 = @(-50388, -73758, 23369, 210474, 137434, 73040, 260862, 211192, 49671);
    18 = @(-48262, -71377, 23115, 208776, 135806, 72970, 257037, 207183, 49854);
    19 = @(-42720, -67038, 24318, 210157, 136299, 73857, 252877, 203338, 49539);
    20 = @(-40596, -64945, 24349, 209739, 136127, 73611, 250335, 201072, 49263);
    21 = @(-45421, -69191, 23769, 209883, 136699, 73184, 255304, 205889, 49415);
    23 = @(-676684, -922026, 245342, 2134441, 1428798, 705643, 2811125, 2350825, 460301);
    24 = @(-242122, -351527, 109405, 893582, 584220, 309362, 1135704, 935747, 199958);
    25 = @(-45452, -67839, 22387, 205091, 135567, 69524, 250543, 203406, 47137);
    26 = @(-41639, -63702, 22063, 204819, 135701, 69118, 246458, 199403, 47055);
    27 = @(-47157, -68718, 21561, 187490, 126875, 60615, 234647, 195594, 39053);
    28 = @(-52959, -74616, 21657, 150074, 95025, 55049, 203033, 169641, 33392);
    29 = @(-54915, -76652, 21736, 146108, 91051, 55056, 201023, 167703, 33320);
    30 = @(-50675, -72004, 21329, 158805, 103702, 55103, 209480, 175706, 33774);
    31 = @(-60743, -80792, 20049, 170908, 115880, 55028, 231651, 196672, 34979);
    32 = @(-63733, -82966, 19233, 174287, 118981, 55306, 238020, 201947, 36073);
    33 = @(-62625, -81613, 18987, 178063, 121965, 56099, 240689, 203577, 37111);
    34 = @(-63678, -82127, 18450, 182732, 125761, 56971, 246410, 207888, 38522);
    35 = @(-67307, -86227, 18920, 185186, 126789, 58397, 252494, 213017, 39477);
    36 = @(-65802, -84770, 18969, 190877, 131500, 59377, 256678, 216270, 40408);
    38 = @(-353070, -444504, 91434, 995172, 699196, 295976, 1348242, 1143700, 204542);
    39 = @(-67092, -86444, 19352, 193221, 134486, 58735, 260313, 220930, 39383);
    40 = @(-70643, -89155, 18513, 188561, 130436, 58125, 259203, 219591, 39612);
    41 = @(-75025, -92859, 17834, 202669, 143658, 59011, 277693, 236516, 41177);
    42 = @(-69071, -86873, 17802, 204704, 145088, 59615, 273775, 231961, 41814);
    43 = @(-71240, -89174, 17934, 206018, 145528, 60490, 277259, 234702, 42556);
    52 = @(-68899, -86680, 17781, 204992, 145288, 59704, 273891, 231968, 41923)
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    $col = 2   # column B
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}

Write-Output "Exhibit 1 refreshed through April; recession-shading rows updated."